$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '65.590.81'
    'E2' = '  -1.47%  '
    'D3' = '3.524.96'
    'E4' = '  +0.02%  '
    'D5' = '599.62'
    'E5' = '  -1.16%  '
    'D6' = '143.27'
    'E6' = '  -1.21%  '
    'D7' = '3.524.08'
    'E7' = '  -1.04%  '
    'D8' = '0.999'
    'E8' = '  -0.15%  '
    'D9' = '0.504'
    'E9' = '  +2.27%  '
    'E10' = '  -1.04%  '
    'E11' = '  -2.97%  '
    'D12' = '0.403'
    'E12' = '  -2.94%  '
    'D13' = '4.127.38'
    'E13' = '  -0.90%  '
    'D14' = '0.0000198'
    'E14' = '  -4.25%  '
    'D15' = '28.59'
    'E15' = '  -4.79%  '
    'D16' = '3.527.62'
    'E16' = '  -0.09%  '
    'E17' = '  +1.26%  '
    'D18' = '65.596.84'
    'E18' = '  -1.46%  '
    'D19' = '10.94'
    'E19' = '  -4.81%  '
    'E20' = '  -0.39%  '
    'D21' = '14.31'
    'E21' = '  -4.69%  '
    'D22' = '414.88'
    'E22' = '  -3.93%  '
    'D23' = '0.598'
    'E23' = '  -2.96%  '
    'D24' = '77.33'
    'E24' = '  -2.44%  '
    'D25' = '3.667.80'
    'E25' = '  -0.89%  '
    'D26' = '0.999'
    'E26' = '  -0.09%  '
    'E27' = '  -3.19%  '
    'E28' = '  -2.87%  '
    'D29' = '8.90'
    'E29' = '  -3.02%  '
    'D30' = '7.69'
    'E30' = '  -4.27%  '
    'E31' = '  +0.06%  '
    'D32' = '3.526.39'
    'E32' = '  -0.81%  '
    'E33' = '  -2.09%  '
    'D34' = '24.28'
    'E34' = '  -4.15%  '
    'D36' = '7.53'
    'E36' = '  -3.87%  '
    'E37' = '  -11.68%  '
    'D38' = '174.22'
    'E38' = '  +0.47%  '
    'D39' = '5.28'
    'E39' = '  -6.09%  '
    'D41' = '0.0818'
    'E41' = '  -3.77%  '
    'D42' = '5.07'
    'E42' = '  -2.50%  '
    'E43' = '  -3.66%  '
    'D44' = '45.26'
    'E44' = '  -2.03%  '
    'D45' = '1.75'
    'E45' = '  -9.52%  '
    'E46' = '  +0.06%  '
    'D47' = '2.36'
    'E47' = '  -6.16%  '
    'D48' = '7.07'
    'E48' = '  -1.59%  '
    'D49' = '22.67'
    'E49' = '  -2.71%  '
    'D50' = '1.08'
    'E50' = '  -9.10%  '
    'D51' = '22.62'
    'E51' = '  -9.89%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$addr]
    $cell.Style = 'Normal'
}